# Weekly data refresh: the rows of price observations (rows 2-24) were
# re-shuffled into a new order. For every destination row, the values of
# columns D (Fecha), I (Calidad), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), O (Origen) and
# P (Precio $/Kg) now come from a different source row, per the mapping
# below (destination row -> source row, both referring to the *original*
# row contents).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    2 = 10
    3 = 9
    4 = 12
    5 = 8
    6 = 14
    7 = 4
    8 = 7
    9 = 22
    10 = 3
    11 = 20
    12 = 23
    13 = 16
    14 = 17
    15 = 24
    16 = 15
    17 = 19
    18 = 5
    19 = 13
    20 = 18
    21 = 11
    22 = 2
    23 = 6
    24 = 21
}

$colD = 4
$colI = 9
$colJ = 10
$colK = 11
$colL = 12
$colM = 13
$colO = 15
$colP = 16

# --- Pass 1: snapshot the original values of every row we may need to
# read from, before any cell gets overwritten. ---
$origD = @{}
$origI = @{}
$origJ = @{}
$origK = @{}
$origL = @{}
$origM = @{}
$origO = @{}
$origP = @{}

for ($r = 2; $r -le 24; $r++) {
    $origD[$r] = [double]$ws.Cells.Item($r, $colD).Value2
    $origI[$r] = $ws.Cells.Item($r, $colI).Text
    $origJ[$r] = [double]$ws.Cells.Item($r, $colJ).Value2
    $origK[$r] = [double]$ws.Cells.Item($r, $colK).Value2
    $origL[$r] = [double]$ws.Cells.Item($r, $colL).Value2
    $origM[$r] = [double]$ws.Cells.Item($r, $colM).Value2
    $origO[$r] = $ws.Cells.Item($r, $colO).Text
    $origP[$r] = [double]$ws.Cells.Item($r, $colP).Value2
}

# --- Pass 2: write the shuffled values back out using the snapshot. ---
foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]

    $ws.Cells.Item($destRow, $colD).Value = $origD[$srcRow]
    $ws.Cells.Item($destRow, $colI).Value = $origI[$srcRow]
    $ws.Cells.Item($destRow, $colJ).Value = $origJ[$srcRow]
    $ws.Cells.Item($destRow, $colK).Value = $origK[$srcRow]
    $ws.Cells.Item($destRow, $colL).Value = $origL[$srcRow]
    $ws.Cells.Item($destRow, $colM).Value = $origM[$srcRow]
    $ws.Cells.Item($destRow, $colO).Value = $origO[$srcRow]
    $ws.Cells.Item($destRow, $colP).Value = $origP[$srcRow]
}

Write-Host "Row reshuffle applied."
